# First take of ExportCheckPointAccessRule
# Adds a new "Access Layers" worksheet (after "Security Zones") documenting
# the Check Point Access Layer import/export columns, mirroring the layout
# used by the other "Import.xlsx" reference sheets.

$wb = $excel.ActiveWorkbook

# Add the new sheet immediately after the current last sheet ("Security Zones")
# so it lands at the end of the tab strip and becomes the active tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Access Layers"

# Header row (bold, like every other sheet in this workbook)
$headers = @(
    "Name",
    "AddDefaultRule",
    "ApplicationsAndUrlFiltering",
    "ContentAwareness",
    "DetectUsingXForwardFor",
    "Firewall",
    "MobileAccess",
    "Shared",
    "Color",
    "Comments",
    "Tags"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
}

# Row 2 - minimal example row
$ws.Cells.Item(2, 1).Value = "TestAccessLayer1"
$ws.Cells.Item(2, 9).Value = "Red"

# Row 3 - fully populated example row
$ws.Cells.Item(3, 1).Value = "TestAccessLayer2"
$ws.Cells.Item(3, 2).Value = $false
$ws.Cells.Item(3, 3).Value = $true
$ws.Cells.Item(3, 4).Value = $true
$ws.Cells.Item(3, 5).Value = $true
$ws.Cells.Item(3, 6).Value = $true
$ws.Cells.Item(3, 8).Value = $true
$ws.Cells.Item(3, 9).Value = "Green"

# Column widths matching the other sheets in the workbook
$ws.Range("A1:H1").ColumnWidth = 21.140625
$ws.Range("I1").ColumnWidth = 13.5703125
$ws.Range("J1").ColumnWidth = 40.85546875
$ws.Range("K1").ColumnWidth = 28.140625

# Leave the selection on E3, matching the authored file
[void]$ws.Range("E3").Select()
